$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 currently holds the old "[XRAY - nbparsingservice, nbexecworker,
# jupyter-notebook-userpod] Make base image tag more specific" bug (ID 2156702).
# That bug got re-titled and re-filed at the bottom of the list, so remove it
# from its old spot here; everything below shifts up by one row.
$ws.Rows.Item(24).Delete()

# Append the two new / updated rows at the bottom of the table (now rows 31 and 32).
# The re-filed "nbparsingservice..." bug title is entered first (matching the
# original authoring order), followed by the row 31 fields and the brand new
# "system.drawing.common" bug.
$ws.Range("C32").Value = "[XRAY] nbparsingservice, nbexecworker, jupyter-notebook-userpod - Make base image tag more specific "
$ws.Range("A32").Value = 2156702
$ws.Range("B32").Value = "Bug"
$ws.Range("D32").Value = "Closed"

$ws.Range("A31").Value = 2203752
$ws.Range("B31").Value = "Bug"
$ws.Range("C31").Value = "[XRAY] ""system.drawing.common"" package on .NET builds"
$ws.Range("D31").Value = "Closed"

# New rows pick up a date-style number format in column G (left blank otherwise).
$ws.Range("G31").NumberFormat = "m/d/yy h:mm"
$ws.Range("G32").NumberFormat = "m/d/yy h:mm"

# Match the saved selection state from the edit.
$ws.Range("F31:G32").Select()
